$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows right before the current row 435, shifting the
# existing rows 435:461 down to become 437:463 (same content, unchanged).
$ws.Rows("435:436").Insert()

# New row 435
$ws.Range("A435").Value = 9
$ws.Range("B435").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C435").Value = "Metropolitana"
$ws.Range("D435").Value = 44931
$ws.Range("E435").Value = 13
$ws.Range("F435").Value = 100112052
$ws.Range("G435").Value = "Albahaca"
$ws.Range("H435").Value = "Sin especificar"
$ws.Range("I435").Value = "Primera"
$ws.Range("J435").Value = 430
$ws.Range("K435").Value = 4000
$ws.Range("L435").Value = 5000
$ws.Range("M435").Value = 4500
$ws.Range("N435").Value = "$/atado"
$ws.Range("O435").Value = "Provincia de Chacabuco"
$ws.Range("P435").Value = 4500
$ws.Range("Q435").Value = 1
$ws.Range("R435").Value = "Hortaliza"

# New row 436
$ws.Range("A436").Value = 9
$ws.Range("B436").Value = "Vega Central Mapocho de Santiago"
$ws.Range("C436").Value = "Metropolitana"
$ws.Range("D436").Value = 44931
$ws.Range("E436").Value = 13
$ws.Range("F436").Value = 100112052
$ws.Range("G436").Value = "Albahaca"
$ws.Range("H436").Value = "Sin especificar"
$ws.Range("I436").Value = "Primera"
$ws.Range("J436").Value = 340
$ws.Range("K436").Value = 4000
$ws.Range("L436").Value = 5000
$ws.Range("M436").Value = 4500
$ws.Range("N436").Value = "$/atado"
$ws.Range("O436").Value = "Región de O'Higgins"
$ws.Range("P436").Value = 4500
$ws.Range("Q436").Value = 1
$ws.Range("R436").Value = "Hortaliza"

# Keep the date-formatted style consistent with the rest of column D
$ws.Range("D435").NumberFormat = "YYYY-MM-DD HH:MM:SS"
$ws.Range("D436").NumberFormat = "YYYY-MM-DD HH:MM:SS"
